$wb = $excel.ActiveWorkbook

# --- Identify existing sheets -------------------------------------------------
# Sheet 1: "总计" (totals) summary sheet
# Sheet 2: "2022-Q2" - holds the quarterly fund-holding detail (will be renamed
#          to "2022-Q4" and given new data); a duplicate of its current
#          (old "2022-Q2") content is preserved on a brand-new sheet.
$wsTotal = $wb.Worksheets.Item(1)
$wsQuarter = $wb.Worksheets.Item(2)

# --- Step 1: preserve the existing "2022-Q2" sheet on a new tab -------------
# Insert a fresh sheet right after the quarter sheet while it still holds the
# original ("2022-Q2") data, then copy that data + formatting across.
$wsPreserved = $wb.Worksheets.Add($null, $wsQuarter)
$wsQuarter.Range("A1:H6").Copy($wsPreserved.Range("A1:H6"))

# Rename: the original tab becomes Q4 first so the "2022-Q2" name is free for
# the newly-added tab to take.
$wsQuarter.Name = "2022-Q4"
$wsPreserved.Name = "2022-Q2"

# --- Step 2: replace sheet 2's contents with the new 2022-Q4 data -----------
# Force text formatting on the numeric-looking detail columns (fund codes and
# percentage/price text) so leading zeros / trailing zeros survive, matching
# the source data which stores these as text.
$wsQuarter.Range("B2:G6").NumberFormat = "@"

$wsQuarter.Cells.Item(2, 2).Value = "000593"
$wsQuarter.Cells.Item(2, 3).Value = "易方达标普全球高端消费品指数增强（QDII）美元现汇"
$wsQuarter.Cells.Item(2, 4).Value = "2.30"
$wsQuarter.Cells.Item(2, 5).Value = "93.71"
$wsQuarter.Cells.Item(2, 6).Value = "3.61"
$wsQuarter.Cells.Item(2, 7).Value = "0.0830"
$wsQuarter.Cells.Item(2, 8).Value = 10

$wsQuarter.Cells.Item(3, 2).Value = "005676"
$wsQuarter.Cells.Item(3, 3).Value = "易方达标普全球高端消费品指数增强C（QDII）人民币"
$wsQuarter.Cells.Item(3, 4).Value = "2.30"
$wsQuarter.Cells.Item(3, 5).Value = "93.71"
$wsQuarter.Cells.Item(3, 6).Value = "3.61"
$wsQuarter.Cells.Item(3, 7).Value = "0.0830"
$wsQuarter.Cells.Item(3, 8).Value = 10

$wsQuarter.Cells.Item(4, 2).Value = "118002"
$wsQuarter.Cells.Item(4, 3).Value = "易方达标普全球高端消费品指数增强A（QDII）人民币"
$wsQuarter.Cells.Item(4, 4).Value = "2.30"
$wsQuarter.Cells.Item(4, 5).Value = "93.71"
$wsQuarter.Cells.Item(4, 6).Value = "3.61"
$wsQuarter.Cells.Item(4, 7).Value = "0.0830"
$wsQuarter.Cells.Item(4, 8).Value = 10

$wsQuarter.Cells.Item(5, 2).Value = "010343"
$wsQuarter.Cells.Item(5, 3).Value = "华宝英国富时100指数A"
$wsQuarter.Cells.Item(5, 4).Value = "0.14"
$wsQuarter.Cells.Item(5, 5).Value = "94.75"
$wsQuarter.Cells.Item(5, 6).Value = "4.10"
$wsQuarter.Cells.Item(5, 7).Value = "0.0057"
$wsQuarter.Cells.Item(5, 8).Value = 6

$wsQuarter.Cells.Item(6, 2).Value = "010344"
$wsQuarter.Cells.Item(6, 3).Value = "华宝英国富时100指数C"
$wsQuarter.Cells.Item(6, 4).Value = "0.08"
$wsQuarter.Cells.Item(6, 5).Value = "94.75"
$wsQuarter.Cells.Item(6, 6).Value = "4.10"
$wsQuarter.Cells.Item(6, 7).Value = "0.0033"
$wsQuarter.Cells.Item(6, 8).Value = 6

# Re-apply the bold/centered header style (matching the "总计" sheet's style)
# to the header row and the row-index column now that the values were rewritten.
$wsTotal.Range("B1").Copy()
$wsQuarter.Range("B1:H1").PasteSpecial(-4122)

$wsTotal.Range("A2").Copy()
$wsQuarter.Range("A2:A6").PasteSpecial(-4122)

# Match the "总计" sheet's page margins (0.75in/1in/0.5in in points).
$qps = $wsQuarter.PageSetup
$qps.LeftMargin = 54
$qps.RightMargin = 54
$qps.TopMargin = 72
$qps.BottomMargin = 72
$qps.HeaderMargin = 36
$qps.FooterMargin = 36

# --- Step 3: update the "总计" summary sheet --------------------------------
# Existing row 2 now reports the Q4 numbers; a new row 3 records the old Q2
# figures that used to live in row 2.
$wsTotal.Cells.Item(2, 2).Value = "2022-Q4"
$wsTotal.Cells.Item(2, 4).Value = 0.26

$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = "2022-Q2"
$wsTotal.Cells.Item(3, 3).Value = 5
$wsTotal.Cells.Item(3, 4).Value = 0.23

$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)
